$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 1
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -2
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = -7
$ws.Range("F18").Value = 4
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = 4
